$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.495724333333333
$ws.Range("H2").Value = 25.487173
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.184802666666667
$ws.Range("N2").Value = 6.554408
$ws.Range("O2").Value = 0.5420193487373902
$ws.Range("P2").Value = 0.5420193487373902
$ws.Range("Q2").Value = 18.56148117873155
$ws.Range("R2").Value = 167.053330608584
$ws.Range("S2").Value = 0.5420193487373902
$ws.Range("T2").Value = 0.5420193487373902

# Row 3 updates
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.495724333333333
$ws.Range("H3").Value = 25.487173
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.846054666666667
$ws.Range("N3").Value = 5.538164
$ws.Range("O3").Value = 0.4579806512626098
$ws.Range("P3").Value = 0.4579806512626098
$ws.Range("Q3").Value = 15.68357155226355
$ws.Range("R3").Value = 141.152143970372
$ws.Range("S3").Value = 0.4579806512626098
$ws.Range("T3").Value = 0.4579806512626098
